$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 14917472
$ws.Range("I19").Value = 13416949
$ws.Range("J19").Value = 16668083
$ws.Range("K19").Value = 13416949
$ws.Range("L19").Value = 16668083
$ws.Range("M19").Value = -13416774
$ws.Range("N19").Value = -16668433

$ws.Range("H28").Value = 887.3570999999999
$ws.Range("I28").Value = 401
$ws.Range("J28").Value = 2103.25
$ws.Range("K28").Value = 401
$ws.Range("L28").Value = 2103.25
$ws.Range("M28").Value = 84
$ws.Range("N28").Value = -3073.25

$ws.Range("H132").Value = 2267.1765
$ws.Range("I132").Value = 2236.9656
$ws.Range("J132").Value = 2442.4
$ws.Range("K132").Value = 6710.8968
$ws.Range("L132").Value = 7327.200000000001
$ws.Range("M132").Value = -4180.8968
$ws.Range("N132").Value = -12387.2

$ws.Range("H138").Value = 3262.9429
$ws.Range("I138").Value = 767.6111
$ws.Range("J138").Value = 5905.0586
$ws.Range("K138").Value = 2302.8333
$ws.Range("L138").Value = 17715.1758
$ws.Range("M138").Value = 2837.1667
$ws.Range("N138").Value = -27995.1758

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4022.0476
$ws.Range("I61").Value = 3023.8333
$ws.Range("J61").Value = 5353
$ws.Range("K61").Value = 3023.8333
$ws.Range("L61").Value = 5353
$ws.Range("M61").Value = -2811.8333
$ws.Range("N61").Value = -5777

$ws.Range("H88").Value = 2249.5
$ws.Range("I88").Value = 3500
$ws.Range("J88").Value = 1999.4
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 1999.4
$ws.Range("M88").Value = -3094
$ws.Range("N88").Value = -2811.4

$ws.Range("H91").Value = 2249.5
$ws.Range("I91").Value = 3500
$ws.Range("J91").Value = 1999.4
$ws.Range("K91").Value = 3500
$ws.Range("L91").Value = 1999.4
$ws.Range("M91").Value = -2096
$ws.Range("N91").Value = -4807.4

$ws.Range("H122").Value = 1765332
$ws.Range("I122").Value = 2470688.8
$ws.Range("J122").Value = 1940.1666
$ws.Range("K122").Value = 7412066.399999999
$ws.Range("L122").Value = 5820.4998
$ws.Range("M122").Value = -7409616.399999999
$ws.Range("N122").Value = -10720.4998

$ws.Range("H136").Value = 4022.0476
$ws.Range("I136").Value = 3023.8333
$ws.Range("J136").Value = 5353
$ws.Range("K136").Value = 9071.499899999999
$ws.Range("L136").Value = 16059
$ws.Range("M136").Value = -6521.499899999999
$ws.Range("N136").Value = -21159

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3764.1428
$ws.Range("I86").Value = 3469.8
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 3469.8
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -2346.8
$ws.Range("N86").Value = -6746

$ws.Range("H89").Value = 3764.1428
$ws.Range("I89").Value = 3469.8
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 17349
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -11733
$ws.Range("N89").Value = -33732

$ws.Range("H99").Value = 4204.4
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4204.4
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = 4204.4
$ws.Range("N99").Value = -7200.4
$ws.Range("L99").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("N126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("M126").ClearContents()

$ws.Range("H134").Value = 4073.762
$ws.Range("I134").Value = 1138.5
$ws.Range("J134").Value = 7987.4443
$ws.Range("K134").Value = 3415.5
$ws.Range("L134").Value = 23962.3329
$ws.Range("M134").Value = -880.5
$ws.Range("N134").Value = -29032.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1379840
$ws.Range("I113").Value = 2463508.5
$ws.Range("J113").Value = 625.4545000000001
$ws.Range("K113").Value = 7390525.5
$ws.Range("L113").Value = 1876.3635
$ws.Range("M113").Value = -7388355.5
$ws.Range("N113").Value = -6216.3635

$ws.Range("H131").Value = 850.6458
$ws.Range("I131").Value = 461.66666
$ws.Range("J131").Value = 980.30554
$ws.Range("K131").Value = 1384.99998
$ws.Range("L131").Value = 2940.91662
$ws.Range("M131").Value = 3655.00002
$ws.Range("N131").Value = -13020.91662

$ws.Range("H132").Value = 632885.75
$ws.Range("I132").Value = 1463
$ws.Range("J132").Value = 1444715
$ws.Range("K132").Value = 13167
$ws.Range("L132").Value = 13002435
$ws.Range("M132").Value = -10637
$ws.Range("N132").Value = -13007495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2694.25
$ws.Range("I102").Value = 2736.2856
$ws.Range("J102").Value = 2400
$ws.Range("K102").Value = 2736.2856
$ws.Range("L102").Value = 2400
$ws.Range("M102").Value = -1114.2856
$ws.Range("N102").Value = -5644

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2413.182
$ws.Range("I136").Value = 1100.25
$ws.Range("K136").Value = 3300.75
$ws.Range("M136").Value = -750.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("N48").Value = 0
$ws.Range("L48").ClearContents()

$ws.Range("H122").Value = 2876.2104
$ws.Range("I122").Value = 2151.2727
$ws.Range("K122").Value = 6453.8181
$ws.Range("M122").Value = -4003.8181

$ws.Range("H132").Value = 1235.8085
$ws.Range("I132").Value = 903.3226
$ws.Range("J132").Value = 1880
$ws.Range("K132").Value = 2709.9678
$ws.Range("L132").Value = 5640
$ws.Range("M132").Value = -179.9677999999999
$ws.Range("N132").Value = -10700

$ws.Range("H135").Value = 39000
$ws.Range("J135").Value = 39000
$ws.Range("L135").Value = 39000
$ws.Range("N135").Value = -49140

$ws.Range("H136").Value = 1346.6
$ws.Range("I136").Value = 1084.5
$ws.Range("K136").Value = 3253.5
$ws.Range("M136").Value = -703.5

$ws.Range("H140").Value = 49613
$ws.Range("J140").Value = 49613
$ws.Range("L140").Value = 49613
$ws.Range("N140").Value = -59973
